$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Job")

$ws.Range("G7").Value = "52000104"
$ws.Range("G8").Value = "52000023"
$ws.Range("G9").Value = "52000065,53000099"
$ws.Range("G10").Value = "53000014,53000018"
$ws.Range("G11").Value = "53000015,53000060"
$ws.Range("G12").Value = "53000043"
$ws.Range("H12").Value = "22032006"
$ws.Range("G13").Value = "53000085"
$ws.Range("H13").Value = "22032005"
$ws.Range("G14").Value = "52000038,53000081"

$ws.Activate()
$ws.Range("H13").Select()
